# Slide 11 ("prez_healhub" final slide) edits:
#  1. Move the background picture ("Рисунок 6") slightly.
#  2. Collapse the "VK TG Git - " run-triplet down to a single "- " run
#     in each of the 3 paragraphs of the "Объект 2" placeholder.
#  3. Move/resize-offset the translucent overlay rectangle ("Прямоугольник 9").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)

# --- 1. Picture offset: (9053,0) -> (4527,1) EMU -------------------------
$picture = $s.Shapes.Item(1)
$picture.Left = (4527 + 0.4) / 12700
$picture.Top  = 1.6 / 12700

# --- 2. Text cleanup: "VK TG Git - " -> "- " (3 occurrences) ------------
$bio = $s.Shapes.Item(4)
$bio.TextFrame.TextRange.Replace("VK TG Git - ", "- ")

# --- 3. Overlay rectangle offset: (9354,3721105) -> (-51815,3764761) EMU -
$overlay = $s.Shapes.Item(5)
$overlay.Left = (-51815 - 0.4) / 12700
$overlay.Top  = (3764761 + 0.4) / 12700
